$d = $word.ActiveDocument

$d.Content.Find.Execute("114×9=", $true, $false, $false, $false, $false, $true, 1, $false, "971×7=", 2)
$d.Content.Find.Execute("345×6=", $true, $false, $false, $false, $false, $true, 1, $false, "536×6=", 2)
$d.Content.Find.Execute("725×7=", $true, $false, $false, $false, $false, $true, 1, $false, "433×8=", 2)
$d.Content.Find.Execute("961×3=", $true, $false, $false, $false, $false, $true, 1, $false, "257×8=", 2)
$d.Content.Find.Execute("205×8=", $true, $false, $false, $false, $false, $true, 1, $false, "985×5=", 2)
$d.Content.Find.Execute("824×2=", $true, $false, $false, $false, $false, $true, 1, $false, "676×8=", 2)
$d.Content.Find.Execute("217×2=", $true, $false, $false, $false, $false, $true, 1, $false, "329×2=", 2)
$d.Content.Find.Execute("343×8=", $true, $false, $false, $false, $false, $true, 1, $false, "347×6=", 2)
$d.Content.Find.Execute("810×5=", $true, $false, $false, $false, $false, $true, 1, $false, "197×8=", 2)
$d.Content.Find.Execute("640×7=", $true, $false, $false, $false, $false, $true, 1, $false, "918×9=", 2)
$d.Content.Find.Execute("285×4=", $true, $false, $false, $false, $false, $true, 1, $false, "291×7=", 2)
$d.Content.Find.Execute("188×5=", $true, $false, $false, $false, $false, $true, 1, $false, "346×2=", 2)
$d.Content.Find.Execute("412×3=", $true, $false, $false, $false, $false, $true, 1, $false, "475×9=", 2)
$d.Content.Find.Execute("986×8=", $true, $false, $false, $false, $false, $true, 1, $false, "380×3=", 2)
$d.Content.Find.Execute("683×5=", $true, $false, $false, $false, $false, $true, 1, $false, "310×5=", 2)
$d.Content.Find.Execute("682×6=", $true, $false, $false, $false, $false, $true, 1, $false, "713×5=", 2)
$d.Content.Find.Execute("361×2=", $true, $false, $false, $false, $false, $true, 1, $false, "773×4=", 2)
$d.Content.Find.Execute("543×2=", $true, $false, $false, $false, $false, $true, 1, $false, "408×6=", 2)
$d.Content.Find.Execute("820×9=", $true, $false, $false, $false, $false, $true, 1, $false, "242×4=", 2)
$d.Content.Find.Execute("691×9=", $true, $false, $false, $false, $false, $true, 1, $false, "574×4=", 2)
$d.Content.Find.Execute("822×4=", $true, $false, $false, $false, $false, $true, 1, $false, "433×8=", 2)
$d.Content.Find.Execute("838×3=", $true, $false, $false, $false, $false, $true, 1, $false, "370×9=", 2)
$d.Content.Find.Execute("163×3=", $true, $false, $false, $false, $false, $true, 1, $false, "436×5=", 2)
$d.Content.Find.Execute("969×2=", $true, $false, $false, $false, $false, $true, 1, $false, "136×6=", 2)
$d.Content.Find.Execute("991×6=", $true, $false, $false, $false, $false, $true, 1, $false, "333×9=", 2)
